$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new BOM rows (10-41) ---
# Column A (row numbers) written first row by row so numeric cells are in place
$ws.Range("A10").Value = 1
$ws.Range("A11").Value = 2
$ws.Range("A12").Value = 3
$ws.Range("A13").Value = 4
$ws.Range("A14").Value = 5
$ws.Range("A15").Value = 6
$ws.Range("A16").Value = 7
$ws.Range("A17").Value = 8
$ws.Range("A18").Value = 9
$ws.Range("A19").Value = 10
$ws.Range("A20").Value = 11
$ws.Range("A21").Value = 12
$ws.Range("A22").Value = 13
$ws.Range("A23").Value = 14
$ws.Range("A24").Value = 15
$ws.Range("A25").Value = 16
$ws.Range("A26").Value = 17
$ws.Range("A27").Value = 18
$ws.Range("A28").Value = 19
$ws.Range("A29").Value = 20
$ws.Range("A30").Value = 21
$ws.Range("A31").Value = 22
$ws.Range("A32").Value = 23
$ws.Range("A33").Value = 24
$ws.Range("A34").Value = 25
$ws.Range("A35").Value = 26
$ws.Range("A36").Value = 27
$ws.Range("A37").Value = 28
$ws.Range("A38").Value = 29
$ws.Range("A39").Value = 30
$ws.Range("A40").Value = 31
$ws.Range("A41").Value = 32

# String cells written in the exact order needed to reproduce the shared-string table
$ws.Range("B10").Value = "nRF52832 module"
$ws.Range("C10").Value = "U1"
$ws.Range("D10").Value = "Bluetooth module"
$ws.Range("C15").Value = "L1"
$ws.Range("B15").Value = "10uH Inductor"
$ws.Range("B13").Value = "10uF capacitor"
$ws.Range("B16").Value = "unknown capcitor"
$ws.Range("C16").Value = "C5, C7"
$ws.Range("B17").Value = "crystal"
$ws.Range("C17").Value = "Y1"
$ws.Range("B18").Value = "1k resistor"
$ws.Range("B19").Value = "red LED"
$ws.Range("B20").Value = "blue LED"
$ws.Range("C20").Value = "LED2"
$ws.Range("B21").Value = "47p capcitor"
$ws.Range("C21").Value = "C10, C11"
$ws.Range("B22").Value = "27 resistor"
$ws.Range("C22").Value = "R11, R12"
$ws.Range("B23").Value = "100nF capacitor"
$ws.Range("B24").Value = "FT231XS"
$ws.Range("C24").Value = "U2"
$ws.Range("B25").Value = "JUMPER-SMT_2_NO_SILK"
$ws.Range("C25").Value = "JP1"
$ws.Range("C18").Value = "R7, R8, R14, R15"
$ws.Range("C19").Value = "LED1, LED3"
$ws.Range("B26").Value = "yellow LED"
$ws.Range("C26").Value = "LED4"
$ws.Range("B27").Value = "0 resistor"
$ws.Range("C27").Value = "R16, R17"
$ws.Range("B28").Value = "OLED "
$ws.Range("C28").Value = "U3"
$ws.Range("B29").Value = "2.2uF capacitor"
$ws.Range("C29").Value = "C16, C14, C17, C18, C15"
$ws.Range("C13").Value = "C4, C12"
$ws.Range("B30").Value = "390k resistor"
$ws.Range("C30").Value = "R13"
$ws.Range("B31").Value = "CONN_02"
$ws.Range("C31").Value = "J1"
$ws.Range("B32").Value = "MIC5504"
$ws.Range("C32").Value = "LDO1"
$ws.Range("D33").Value = "USB.MICRO-BIGGERPADS"
$ws.Range("B33").Value = "Micro USB Pad"
$ws.Range("C33").Value = "USB1"
$ws.Range("B34").Value = "MCP73831"
$ws.Range("C34").Value = "CHRG1"
$ws.Range("C14").Value = "C2, C3, C6"
$ws.Range("B35").Value = "Schottky diode"
$ws.Range("D35").Value = "SCHOTTKY-USC"
$ws.Range("B37").Value = "100k resistor"
$ws.Range("B38").Value = "5k Resistor"
$ws.Range("C38").Value = "R4"
$ws.Range("C36").Value = "C1, C8"
$ws.Range("B36").Value = "4.7uF capcitor"
$ws.Range("C37").Value = "R3, R5, R6"
$ws.Range("B39").Value = "CORTEX_DEBUG_PTH"
$ws.Range("C39").Value = "DEBUG1"
$ws.Range("C35").Value = "D1, D2, D3"
$ws.Range("C23").Value = "C13, C9"
$ws.Range("C11").Value = "R1, R2, R10, R9"
$ws.Range("C12").Value = "SW1, RESET"
$ws.Range("B40").Value = "CONN_07-1.27MM"
$ws.Range("C40").Value = "J2"
$ws.Range("B12").Value = "Button Gullwig"
$ws.Range("B41").Value = "Buttons MOM"
$ws.Range("C41").Value = "SW3, SW4"

# Reused existing shared strings (B11 -> "10k resistor", B14 -> "1uF capacitor")
$ws.Range("B11").Value = "10k resistor"
$ws.Range("B14").Value = "1uF capacitor"

# --- Update view state: scroll/selection ---
$ws.Range("C30").Select()
